$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph.
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
# The new paragraph must not keep the Heading1 style of its predecessor.
$metaPara.Range.Style = "Normal"

$metaText = "Meta description" + ": Discover the features of Black Mamba slot and play for free. Read our review to find out about its graphics, RTP, volatility, and bonuses."
$metaRange = $metaPara.Range
$metaBody = $d.Range($metaRange.Start, $metaRange.End - 1)
$metaBody.Text = $metaText

# Bold just the "Meta description" label.
$labelLen = ("Meta description").Length
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + $labelLen)
$labelRange.Bold = 1

# Match the source markup's leading empty run before the visible text
# (every sibling body paragraph in this document starts with one).
$leadXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$leadPoint = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$leadPoint.InsertXML($leadXml)

# ------------------------------------------------------------------
# 2) Drop the trailing bold "Play Black Mamba Slot..." paragraph and
#    turn the remaining italic paragraph into the new image-prompt text.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$beforeLastPara = $d.Paragraphs.Item($count - 1)
$beforeLastPara.Range.Delete()

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$newText = 'Please create a feature image fitting the game "Black Mamba" with the following requirements: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses'
$lastRange = $lastPara.Range
$lastBody = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastBody.Text = $newText
